# Saldo.xlsx update
#
# 1) Add a new account row for PEDRO (conta 004460487, saldo 60000).
# 2) Update ANDREA's (conta 005186167) saldo from 144.42 to 294.42.
#
# The sheet is kept sorted by "Saldo" descending, so both rows land at the
# row position that matches their (new) balance rather than being appended
# at the bottom - same as if a user re-sorted the sheet by column C after
# editing the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Insert PEDRO / 004460487 / 60000 ----------------------------------
# 60000 belongs right above THOMAS (008026942 / 20000, row 3) and below
# BLUEMETRIX (97053.05, row 2) to keep the descending sort order.
$thomasRow = $ws.Columns(1).Find("008026942").Row
$ws.Rows($thomasRow).Insert()

$newRow = $thomasRow
$ws.Cells.Item($newRow, 1).NumberFormat = "@"   # keep leading zeros -> text
$ws.Cells.Item($newRow, 1).Value = "004460487"
$ws.Cells.Item($newRow, 2).Value = "PEDRO"
$ws.Cells.Item($newRow, 3).Value = 60000

# --- 2) Move ANDREA / 005186167 to its new sorted position -----------------
# Remove the existing ANDREA row (currently saldo 144.42) ...
$andreaRow = $ws.Columns(1).Find("005186167").Row
$ws.Rows($andreaRow).Delete()

# ... and re-insert it with the updated saldo (294.42) just above THIAGO
# (004381095 / 283.81) and below GILSON (299.6), which is where it now
# belongs in the descending sort order.
$thiagoRow = $ws.Columns(1).Find("004381095").Row
$ws.Rows($thiagoRow).Insert()

$andreaNewRow = $thiagoRow
$ws.Cells.Item($andreaNewRow, 1).NumberFormat = "@"   # keep leading zeros -> text
$ws.Cells.Item($andreaNewRow, 1).Value = "005186167"
$ws.Cells.Item($andreaNewRow, 2).Value = "ANDREA"
$ws.Cells.Item($andreaNewRow, 3).Value = 294.42
